# Addition of word 'escribe'
# Extends the state-transition matrix on Hoja1 (sheet1) from row 28 down
# through new rows 29-36, filling in "ER" (error state) fill cells and the
# diagonal-ish numbered transition cells, then appends the new PR03/escribe
# pair of cells on row 36 (mirroring PR02/principal on row 28).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row whose formatting (fill s="2") we copy down onto column A of the new rows.
$aStyleSource = $ws.Range("A28")
# Row whose formatting (underline font s="3") we copy onto AK35.
$uStyleSource = $ws.Range("AM18")

# Diagonal-like "next state" cell per row: column letter + numeric value.
$specialCol = @{
    29 = "T";  30 = "D";  31 = "S";  32 = "J";  33 = "C";  34 = "F"
}
$specialVal = @{
    29 = 21;  30 = 22;  31 = 23;  32 = 24;  33 = 25;  34 = 26
}

for ($row = 29; $row -le 34; $row++) {
    # Column A: sequence number, with the gray "s=2" fill style copied down.
    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value = $row - 9
    $aStyleSource.Copy()
    $aCell.PasteSpecial(-4122)

    # Columns B:AO (2..41): default to the "ER" shared string.
    for ($col = 2; $col -le 41; $col++) {
        $ws.Cells.Item($row, $col).Value = "ER"
    }

    # Overwrite the one "next state" cell for this row with its number.
    $col = $specialCol[$row]
    $val = $specialVal[$row]
    $ws.Range($col + $row).Value = $val
}

# --- Row 35 ---------------------------------------------------------------
$aCell = $ws.Cells.Item(35, 1)
$aCell.Value = 26
$aStyleSource.Copy()
$aCell.PasteSpecial(-4122)

for ($col = 2; $col -le 41; $col++) {
    $ws.Cells.Item(35, $col).Value = "ER"
}

# AK35 keeps the "ER" text but picks up the underline font style (s="3").
$ws.Range("AK35").Value = "ER"
$uStyleSource.Copy()
$ws.Range("AK35").PasteSpecial(-4122)

$ws.Range("AL35").Value = 27
$ws.Range("AM35").Value = 27
$ws.Range("AN35").Value = 27

# --- Row 36 ---------------------------------------------------------------
$aCell = $ws.Cells.Item(36, 1)
$aCell.Value = 27
$aStyleSource.Copy()
$aCell.PasteSpecial(-4122)

for ($col = 2; $col -le 41; $col++) {
    $ws.Cells.Item(36, $col).Value = "ER"
}

# New vocabulary entry: PR03 / escribe (mirrors PR02 / principal on row 28).
$ws.Range("AP36").Value = "PR03"
$ws.Range("AQ36").Value = "escribe"

# --- Sheet view: scrolled down, new active selection -----------------------
$view = $ws.Application.ActiveWindow
$view.ScrollRow = 19
$ws.Range("AP41").Select()
